$wb = $excel.ActiveWorkbook

# Sheet: PV Dispatch
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 13.6
$ws.Range("H2").Value = 27.2
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 40.8
$ws.Range("K2").Value = 47.6
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68
$ws.Range("O2").Value = 61.2
$ws.Range("P2").Value = 54.4
$ws.Range("Q2").Value = 47.6
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 20.4
$ws.Range("T2").Value = 13.6
$ws.Range("I3").Value = 27.2
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 68
$ws.Range("N3").Value = 54.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 34
$ws.Range("Q3").Value = 34
$ws.Range("R3").Value = 20.4
$ws.Range("S3").Value = 13.6
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 54.4
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 10.38312417100186

# Sheet: Battery Input
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 64.3
$ws.Range("H2").Value = 149.834272013061
$ws.Range("I2").Value = 2.8
$ws.Range("J2").Value = 1.8
$ws.Range("K2").Value = 21.6
$ws.Range("L2").Value = 33.6
$ws.Range("M2").Value = 37.8
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 30
$ws.Range("P2").Value = 25.8
$ws.Range("Q2").Value = 21.6
$ws.Range("R2").Value = 0.2
$ws.Range("S2").Value = 7.2
$ws.Range("T2").Value = 45.6
$ws.Range("I3").Value = 27.43079277624771
$ws.Range("J3").Value = 40.8
$ws.Range("K3").Value = 54.4
$ws.Range("L3").Value = 61.2
$ws.Range("M3").Value = 44.6
$ws.Range("N3").Value = 28.4
$ws.Range("O3").Value = 47.6
$ws.Range("P3").Value = 5.4
$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 20.4
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 27.2
$ws.Range("L4").Value = 47.6
$ws.Range("M4").Value = 31
$ws.Range("N4").Value = 54.4
$ws.Range("O4").Value = 47.6
$ws.Range("P4").Value = 27.2
$ws.Range("Q4").Value = 10.38312417100186
$ws.Range("R4").Value = 0

# Sheet: State of Charge
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("E2").Value = 133.1313131313138
$ws.Range("G2").Value = 183.657
$ws.Range("H2").Value = 331.9929292929304
$ws.Range("I2").Value = 334.7649292929304
$ws.Range("J2").Value = 336.5469292929304
$ws.Range("K2").Value = 357.9309292929304
$ws.Range("L2").Value = 391.1949292929304
$ws.Range("M2").Value = 428.6169292929304
$ws.Range("N2").Value = 470.1969292929304
$ws.Range("O2").Value = 499.8969292929304
$ws.Range("P2").Value = 525.4389292929304
$ws.Range("Q2").Value = 546.8229292929304
$ws.Range("R2").Value = 547.0209292929304
$ws.Range("S2").Value = 554.1489292929305
$ws.Range("I3").Value = 147.1564848484852
$ws.Range("J3").Value = 187.5484848484852
$ws.Range("K3").Value = 241.4044848484852
$ws.Range("L3").Value = 301.9924848484852
$ws.Range("M3").Value = 346.1464848484852
$ws.Range("N3").Value = 374.2624848484852
$ws.Range("O3").Value = 421.3864848484852
$ws.Range("P3").Value = 426.7324848484852
$ws.Range("Q3").Value = 434.6524848484852
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 146.928
$ws.Range("L4").Value = 194.052
$ws.Range("M4").Value = 224.742
$ws.Range("N4").Value = 278.598
$ws.Range("O4").Value = 325.722
$ws.Range("P4").Value = 352.65
$ws.Range("Q4").Value = 362.9292929292918

# Sheet: Feed in from Type 2
$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("H2").Value = 135.634272013061
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("T2").Value = 34
$ws.Range("I3").Value = 0.2307927762477106
$ws.Range("N3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("S3").Value = 9.6
$ws.Range("K4").Value = 0

# Sheet: Feed in from Type 3
$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

# Sheet: Feed in from Type 4
$ws = $wb.Worksheets.Item("Feed in from Type 4")
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

# Sheet: Costs and Revenues
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 77388.66797673714
$ws.Range("D2").Value = 9992.97670278544
$ws.Range("E2").Value = 1770
$ws.Range("F2").Value = 16326.78438529687

# Sheet: Capacities
$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C3").Value = 68

Write-Output "Applied all edits"